$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for rows 2 through 16
# from serial date 45224 (2023-10-25) to 45233 (2023-11-03).
for ($row = 2; $row -le 16; $row++) {
    $ws.Cells.Item($row, 3).Value = 45233
}
